$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, never letting Excel
# auto-convert numeric-looking strings (e.g. "117.03") into numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "51.761.81"
$ws.Range("E2").Value = "  +6.04%  "

Set-TextValue $ws.Range("D3") "2.750.53"
$ws.Range("E3").Value = "  +4.17%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws.Range("D5") "117.03"
$ws.Range("E5").Value = "  +5.80%  "

Set-TextValue $ws.Range("D6") "332.23"
$ws.Range("E6").Value = "  +2.93%  "

$ws.Range("E7").Value = "  +2.54%  "

$ws.Range("E8").Value = "  -0.03%  "

Set-TextValue $ws.Range("D9") "0.574"
$ws.Range("E9").Value = "  +6.25%  "

Set-TextValue $ws.Range("D10") "41.80"
$ws.Range("E10").Value = "  +5.78%  "

$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("E13").Value = "  +2.88%  "

Set-TextValue $ws.Range("D14") "7.59"
$ws.Range("E14").Value = "  +5.23%  "

Set-TextValue $ws.Range("D15") "3.186.89"
$ws.Range("E15").Value = "  +4.65%  "

Set-TextValue $ws.Range("D16") "2.767.06"
$ws.Range("E16").Value = "  +5.07%  "

Set-TextValue $ws.Range("D17") "0.882"
$ws.Range("E17").Value = "  +1.93%  "

Set-TextValue $ws.Range("D18") "51.744.82"

Set-TextValue $ws.Range("D19") "13.50"
$ws.Range("E19").Value = "  +5.32%  "

Set-TextValue $ws.Range("D20") "3.06"
$ws.Range("E20").Value = "  +6.19%  "

$ws.Range("E21").Value = "  +2.34%  "

Set-TextValue $ws.Range("D22") "0.0₃0961"
$ws.Range("E22").Value = "  +2.13%  "

Set-TextValue $ws.Range("D23") "278.44"
$ws.Range("E23").Value = "  +2.93%  "

Set-TextValue $ws.Range("D24") "69.79"
$ws.Range("E24").Value = "  +0.72%  "

Set-TextValue $ws.Range("D25") "2.64"
$ws.Range("E25").Value = "  +4.41%  "

Set-TextValue $ws.Range("D26") "26.76"
$ws.Range("E26").Value = "  +2.63%  "

Set-TextValue $ws.Range("D27") "4.14"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("E28").Value = "  +0.06%  "

Set-TextValue $ws.Range("D29") "10.24"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("E31").Value = "  +2.10%  "

Set-TextValue $ws.Range("D32") "35.01"
$ws.Range("E32").Value = "  +0.18%  "

Set-TextValue $ws.Range("D33") "50.47"
$ws.Range("E33").Value = "  +2.04%  "

Set-TextValue $ws.Range("D34") "5.57"
$ws.Range("E34").Value = "  +2.68%  "

Set-TextValue $ws.Range("D35") "0.0820"
$ws.Range("E35").Value = "  +2.84%  "

$ws.Range("E36").Value = "  -0.04%  "

Set-TextValue $ws.Range("D37") "19.08"
$ws.Range("E37").Value = "  -0.89%  "

Set-TextValue $ws.Range("D38") "2.08"
$ws.Range("E38").Value = "  +2.32%  "

Set-TextValue $ws.Range("D39") "4.96"
$ws.Range("E39").Value = "  +0.96%  "

Set-TextValue $ws.Range("D40") "3.17"
$ws.Range("E40").Value = "  +0.38%  "

Set-TextValue $ws.Range("D41") "129.93"
$ws.Range("E41").Value = "  +3.71%  "

$ws.Range("E44").Value = "  +2.58%  "

Set-TextValue $ws.Range("D47") "2.118.20"
$ws.Range("E47").Value = "  +2.44%  "

Set-TextValue $ws.Range("D48") "3.33"
$ws.Range("E48").Value = "  +3.47%  "

$ws.Range("E49").Value = "  +2.22%  "

Set-TextValue $ws.Range("D50") "5.57"
$ws.Range("E50").Value = "  +7.94%  "

$ws.Range("E51").Value = "  +0.01%  "

# Rows 42/43 and 45/46 swap Coin/Link identity (ranking reorder),
# each also carrying independently updated Price/Volume values.
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "23.29"
$ws.Range("E42").Value = "  +2.43%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D43") "0.0346"
$ws.Range("E43").Value = "  +10.48%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "2.26"
$ws.Range("E45").Value = "  +5.50%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D46") "2.45"
$ws.Range("E46").Value = "  +15.74%  "
